$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.0738400261005589, 15.88588071331727, 4.571465910586882, 0, "eyescrunching"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.1,  -0.02903273376741133, 19.83350219973457, -15.48431052625258, 1, "jaw"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.3664019806239966,  15.83675581695562, 10.24701591653692,  2, "raisingeyebrows"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.02541153791507557, 21.6287663825487,  16.11835549118687,  3, "movehat"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -1.449764027933893,   3.212056155844523, 2.204190255890288,  4, "movehead")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
